# "adding averages and more checks"
#
# Refresh the Training Dashboard's "PERIOD TO EXPIRE" / "LAST UPDATE"
# columns against a later run date, flip the now-expired LOTO SOP row
# from VALID to NOT VALID (re-styling it like the other invalid row),
# white-out the header row font, and shrink + simplify the Exam
# Dashboard's COMMENTS column now that every exam date checks out.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# Helper: write a literal text value (never let Excel's autodetect turn a
# date-shaped string like "16-Sep-2025" into a real date/serial number).
function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- Training Dashboard ---------------------------------------------

$ws1.Range("H3").Value = 394
Set-TextValue $ws1.Range("I3") "16-Sep-2025"

$ws1.Range("H4").Value = 413
Set-TextValue $ws1.Range("I4") "16-Sep-2025"

$ws1.Range("H5").Value = 413
Set-TextValue $ws1.Range("I5") "16-Sep-2025"

# Row 6 (LOTO SOP) has now expired -> NOT VALID. Pull over the red/pink
# "NOT VALID" formatting already used on row 7 so the whole row matches.
$ws1.Range("A7:K7").Copy()
$ws1.Range("A6:K6").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("H6").Value = 14
Set-TextValue $ws1.Range("I6") "16-Sep-2025"
$ws1.Range("J6").Value = "NOT VALID"

$ws1.Range("H7").Value = -328
Set-TextValue $ws1.Range("I7") "16-Sep-2025"

$ws1.Range("H8").Value = 155
Set-TextValue $ws1.Range("I8") "16-Sep-2025"

# Header row: bold white text on its existing dark-blue fill.
$ws1.Range("A2:K2").Font.Color = 16777215

# --- Exam Dashboard ----------------------------------------------------

$ws2.Columns("E").ColumnWidth = 14.166666666666666

$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
$ws2.Range("E5").Value = "date is valid"
